$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 78; existing rows 78-173 shift down to 80-175.
$ws.Rows("78:79").Insert()

# Populate new row 78 (Ciruela / Fortuna / Primera)
$ws.Cells.Item(78, 1).Value = 5
$ws.Cells.Item(78, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(78, 3).Value = "Maule"
$ws.Cells.Item(78, 4).Value = 44966
$ws.Cells.Item(78, 5).Value = 7
$ws.Cells.Item(78, 6).Value = "Fruta"
$ws.Cells.Item(78, 7).Value = 100103
$ws.Cells.Item(78, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(78, 9).Value = 100103002
$ws.Cells.Item(78, 10).Value = "Ciruela"
$ws.Cells.Item(78, 11).Value = "Fortuna"
$ws.Cells.Item(78, 12).Value = "Primera"
$ws.Cells.Item(78, 13).Value = 400
$ws.Cells.Item(78, 14).Value = 10000
$ws.Cells.Item(78, 15).Value = 10000
$ws.Cells.Item(78, 16).Value = 10000
$ws.Cells.Item(78, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(78, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(78, 19).Value = 556
$ws.Cells.Item(78, 20).Value = 18

# Populate new row 79 (Ciruela / Fortuna / Segunda)
$ws.Cells.Item(79, 1).Value = 5
$ws.Cells.Item(79, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(79, 3).Value = "Maule"
$ws.Cells.Item(79, 4).Value = 44966
$ws.Cells.Item(79, 5).Value = 7
$ws.Cells.Item(79, 6).Value = "Fruta"
$ws.Cells.Item(79, 7).Value = 100103
$ws.Cells.Item(79, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(79, 9).Value = 100103002
$ws.Cells.Item(79, 10).Value = "Ciruela"
$ws.Cells.Item(79, 11).Value = "Fortuna"
$ws.Cells.Item(79, 12).Value = "Segunda"
$ws.Cells.Item(79, 13).Value = 150
$ws.Cells.Item(79, 14).Value = 8000
$ws.Cells.Item(79, 15).Value = 8000
$ws.Cells.Item(79, 16).Value = 8000
$ws.Cells.Item(79, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(79, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(79, 19).Value = 444
$ws.Cells.Item(79, 20).Value = 18
